$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> URL, add new sheet "Result" right after it ---
$wsUrl = $wb.Worksheets.Item(1)
$wsUrl.Name = "URL"
$wsResult = $wb.Worksheets.Add($null, $wsUrl)
$wsResult.Name = "Result"

# --- URL sheet: add the two JIRA links in A2 / A3 ---
$wsUrl.Hyperlinks.Add($wsUrl.Range("A2"), "https://jira.jnj.com/browse/AGQP-293")
$wsUrl.Hyperlinks.Add($wsUrl.Range("A3"), "https://jira.jnj.com/browse/AGQP-294")

# --- Result sheet: header row ---
$headers = @("URL", "Title", "Epic", "Type", "Affected Version", "Fix Version", "Story Point", "Acceptance Criteria", "Description", "Priority", "Approval Workflow", "Asignee/Reporter", "Sprint")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsResult.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Result sheet: data rows (one row per JIRA URL) ---
# Color constants (OLE BGR-ish "RGB()" values Excel expects): green = 0x00FF00, yellow = 0xFFFF00
$green = 65280
$yellow = 65535

$urls = @("https://jira.jnj.com/browse/AGQP-293", "https://jira.jnj.com/browse/AGQP-294")

# Columns B..M: text + fill color to apply after the URL column (A)
$rowValues = @("Passed", "Passed", "Update the type as Story", "Passed", "Passed", "Add the Story Point", "Add the acceptance criteria", "Passed", "Passed", "Passed", "Assignee/Reporter Cant be Same", "Passed")
$rowColors = @($green, $green, $yellow, $green, $green, $yellow, $yellow, $green, $green, $green, $yellow, $green)

for ($r = 0; $r -lt $urls.Length; $r++) {
    $row = $r + 2
    $wsResult.Cells.Item($row, 1).Value = $urls[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $cell = $wsResult.Cells.Item($row, $c + 2)
        $cell.Value = $rowValues[$c]
        $cell.Interior.Color = $rowColors[$c]
    }
}

# --- Restore selection / active sheet to URL, cell I5 (matches target selection) ---
$wsUrl.Activate() | Out-Null
$wsUrl.Range("I5").Select() | Out-Null
